$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing task descriptions (column A) -------------------------
$ws.Range("A11").Value = "Game rules (scoring, keeping track of shots and hitting the water)"
$ws.Range("A13").Value = "Physics "
$ws.Range("A14").Value = "In-game UI for shooting"
$ws.Range("A18").Value = "Make it a game / merge / bug fixing"

# --- New rows for the additional planning items ---------------------------
$ws.Range("A19").Value = "Presentation"
$ws.Range("E19").Style = "Good"

$ws.Range("A20").Value = "Handle being out of bounds"
$ws.Range("E20").Style = "Good"

$ws.Range("A21").Value = "UI and terrain improvements (if time allows)"
$ws.Range("E21").Style = "Good"

# --- Sheet-level view / formatting tweaks ----------------------------------
$ws.Columns.Item(1).ColumnWidth = 56

$excel.ActiveWindow.Zoom = 130

$ws.Range("A21").Select() | Out-Null
